# Daily attendance processing - 2025-10-24 12:39:46
# Normalize the "Recorded By" (column G) value ordering: move the
# attendance-taker's own email address to the front of the recorded-by
# list, ahead of the automated "System"/"system" entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, system, System"
    }
    elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "System, backup@backdoor.com") {
        $cell.Value = "backup@backdoor.com, System"
    }
    elseif ($val -eq "admin@admin.com, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, admin@admin.com"
    }
}
